# BOM.xlsx edit — "pull down bij base van bjt"
# 1) D3 note changes from "Hebben we" to "Rauf"
# 2) D9 picks up the same fill/font formatting already used by the rest of
#    column D (it had no explicit style before)
# 3) D10 gets a new note "Hebben we?" with that same column-D formatting
# 4) A new BOM row 11 is added: RPI header / 2x20 / 1 (pull-down resistor's
#    header connector), formatted like the row above it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D3: "Hebben we" -> "Rauf" ---
$ws.Range("D3").Value = "Rauf"

# --- D9: give it the standard column-D formatting ---
$ws.Range("D7").Copy()
$ws.Range("D9").PasteSpecial(-4122)

# --- D10: new note, standard column-D formatting ---
$ws.Range("D10").Value = "Hebben we?"
$ws.Range("D7").Copy()
$ws.Range("D10").PasteSpecial(-4122)

# --- Row 11: new BOM line ---
$ws.Range("A11").Value = "RPI header"
$ws.Range("B11").Value = "2x20"
$ws.Range("C11").Value = 1

$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# Matches the saved selection in the target workbook
$ws.Range("B11").Select()
